# Auto-generated edit script: updates cached numeric values in the
# per-leve price/profit columns (H..N) across all 8 job sheets to
# reflect a refreshed market-board data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 166682.5  # H6: 11.333333 -> 166682.5
$ws.Cells.Item(6, 9).Value = 166682.5  # I6: 11.333333 -> 166682.5
$ws.Cells.Item(6, 11).Value = 500047.5  # K6: 33.999999 -> 500047.5
$ws.Cells.Item(6, 13).Value = -499935.5  # M6: 78.000001 -> -499935.5
$ws.Cells.Item(33, 8).Value = 325.63635  # H33: 328.27274 -> 325.63635
$ws.Cells.Item(33, 10).Value = 393.75  # J33: 408.25 -> 393.75
$ws.Cells.Item(33, 12).Value = 393.75  # L33: 408.25 -> 393.75
$ws.Cells.Item(33, 14).Value = -851.75  # N33: -866.25 -> -851.75
$ws.Cells.Item(47, 8).Value = 9950  # H47: 7975 -> 9950
$ws.Cells.Item(47, 9).Value = 9950  # I47: 7975 -> 9950
$ws.Cells.Item(47, 11).Value = 9950  # K47: 7975 -> 9950
$ws.Cells.Item(47, 13).Value = -8978  # M47: -7003 -> -8978
$ws.Cells.Item(54, 8).Value = 38749.75  # H54: 40000 -> 38749.75
$ws.Cells.Item(54, 9).Value = 38333  # I54: 40000 -> 38333
$ws.Cells.Item(54, 11).Value = 38333  # K54: 40000 -> 38333
$ws.Cells.Item(54, 13).Value = -37847  # M54: -39514 -> -37847
$ws.Cells.Item(75, 8).Value = 76332.336  # H75: 94499.5 -> 76332.336
$ws.Cells.Item(75, 10).Value = 76332.336  # J75: 94499.5 -> 76332.336
$ws.Cells.Item(75, 12).Value = 76332.336  # L75: 94499.5 -> 76332.336
$ws.Cells.Item(75, 14).Value = -78204.336  # N75: -96371.5 -> -78204.336
$ws.Cells.Item(78, 8).Value = 76332.336  # H78: 94499.5 -> 76332.336
$ws.Cells.Item(78, 10).Value = 76332.336  # J78: 94499.5 -> 76332.336
$ws.Cells.Item(78, 12).Value = 228997.008  # L78: 283498.5 -> 228997.008
$ws.Cells.Item(78, 14).Value = -238357.008  # N78: -292858.5 -> -238357.008
$ws.Cells.Item(125, 8).Value = 636.75  # H125: 516 -> 636.75
$ws.Cells.Item(125, 10).Value = 999  # J125: 0 -> 999
$ws.Cells.Item(125, 12).Value = 8991  # L125: 0 -> 8991
$ws.Cells.Item(125, 14).Value = -13911  # N125: None -> -13911

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 434.15384  # H2: 462.75 -> 434.15384
$ws.Cells.Item(2, 9).Value = 395.33334  # I2: 423 -> 395.33334
$ws.Cells.Item(2, 11).Value = 395.33334  # K2: 423 -> 395.33334
$ws.Cells.Item(2, 13).Value = -282.33334  # M2: -310 -> -282.33334
$ws.Cells.Item(61, 8).Value = 8941  # H61: 9671.5 -> 8941
$ws.Cells.Item(61, 9).Value = 8941  # I61: 9671.5 -> 8941
$ws.Cells.Item(61, 11).Value = 8941  # K61: 9671.5 -> 8941
$ws.Cells.Item(61, 13).Value = -8729  # M61: -9459.5 -> -8729
$ws.Cells.Item(116, 8).Value = 434.15384  # H116: 462.75 -> 434.15384
$ws.Cells.Item(116, 9).Value = 395.33334  # I116: 423 -> 395.33334
$ws.Cells.Item(116, 11).Value = 395.33334  # K116: 423 -> 395.33334
$ws.Cells.Item(116, 13).Value = 1898.66666  # M116: 1871 -> 1898.66666
$ws.Cells.Item(132, 8).Value = 3665.9333  # H132: 2769.7727 -> 3665.9333
$ws.Cells.Item(132, 9).Value = 3912.25  # I132: 2651.0881 -> 3912.25
$ws.Cells.Item(132, 11).Value = 11736.75  # K132: 7953.2643 -> 11736.75
$ws.Cells.Item(132, 13).Value = -9206.75  # M132: -5423.2643 -> -9206.75
$ws.Cells.Item(136, 8).Value = 8941  # H136: 9671.5 -> 8941
$ws.Cells.Item(136, 9).Value = 8941  # I136: 9671.5 -> 8941
$ws.Cells.Item(136, 11).Value = 26823  # K136: 29014.5 -> 26823
$ws.Cells.Item(136, 13).Value = -24273  # M136: -26464.5 -> -24273

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 434.15384  # H3: 462.75 -> 434.15384
$ws.Cells.Item(3, 9).Value = 395.33334  # I3: 423 -> 395.33334
$ws.Cells.Item(3, 11).Value = 395.33334  # K3: 423 -> 395.33334
$ws.Cells.Item(3, 13).Value = -281.33334  # M3: -309 -> -281.33334
$ws.Cells.Item(5, 8).Value = 1495  # H5: 1497.6666 -> 1495
$ws.Cells.Item(5, 9).Value = 1495  # I5: 1497.2 -> 1495
$ws.Cells.Item(5, 10).Value = 0  # J5: 1500 -> 0
$ws.Cells.Item(5, 11).Value = 1495  # K5: 1497.2 -> 1495
$ws.Cells.Item(5, 12).Value = 0  # L5: 1500 -> 0
$ws.Cells.Item(5, 13).Value = $null  # M5: -1384.2 -> None
$ws.Cells.Item(5, 14).Value = -1382  # N5: -1726 -> -1382
$ws.Cells.Item(86, 8).Value = 8306.467000000001  # H86: 8978.308000000001 -> 8306.467000000001
$ws.Cells.Item(86, 9).Value = 5035.6665  # I86: 5267 -> 5035.6665
$ws.Cells.Item(86, 10).Value = 10487  # J86: 11297.875 -> 10487
$ws.Cells.Item(86, 11).Value = 5035.6665  # K86: 5267 -> 5035.6665
$ws.Cells.Item(86, 12).Value = 10487  # L86: 11297.875 -> 10487
$ws.Cells.Item(86, 13).Value = -3912.6665  # M86: -4144 -> -3912.6665
$ws.Cells.Item(86, 14).Value = -12733  # N86: -13543.875 -> -12733
$ws.Cells.Item(89, 8).Value = 8306.467000000001  # H89: 8978.308000000001 -> 8306.467000000001
$ws.Cells.Item(89, 9).Value = 5035.6665  # I89: 5267 -> 5035.6665
$ws.Cells.Item(89, 10).Value = 10487  # J89: 11297.875 -> 10487
$ws.Cells.Item(89, 11).Value = 25178.3325  # K89: 26335 -> 25178.3325
$ws.Cells.Item(89, 12).Value = 52435  # L89: 56489.375 -> 52435
$ws.Cells.Item(89, 13).Value = -19562.3325  # M89: -20719 -> -19562.3325
$ws.Cells.Item(89, 14).Value = -63667  # N89: -67721.375 -> -63667
$ws.Cells.Item(132, 8).Value = 0  # H132: 49999 -> 0
$ws.Cells.Item(132, 10).Value = 0  # J132: 49999 -> 0
$ws.Cells.Item(132, 12).Value = $null  # L132: 49999 -> None
$ws.Cells.Item(132, 14).Value = 0  # N132: -60119 -> 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 4337.3335  # H58: 5750 -> 4337.3335
$ws.Cells.Item(58, 9).Value = 1506  # I58: 1500 -> 1506
$ws.Cells.Item(58, 11).Value = 1506  # K58: 1500 -> 1506
$ws.Cells.Item(58, 13).Value = -1303  # M58: -1297 -> -1303
$ws.Cells.Item(132, 8).Value = 2970.2727  # H132: 2862.8696 -> 2970.2727
$ws.Cells.Item(132, 9).Value = 2933.0667  # I132: 2781 -> 2933.0667
$ws.Cells.Item(132, 11).Value = 8799.2001  # K132: 8343 -> 8799.2001
$ws.Cells.Item(132, 13).Value = -6269.2001  # M132: -5813 -> -6269.2001
$ws.Cells.Item(134, 8).Value = 2772.125  # H134: 2997.125 -> 2772.125
$ws.Cells.Item(134, 9).Value = 2375.8  # I134: 2735.8 -> 2375.8
$ws.Cells.Item(134, 11).Value = 7127.400000000001  # K134: 8207.400000000001 -> 7127.400000000001
$ws.Cells.Item(134, 13).Value = -4592.400000000001  # M134: -5672.400000000001 -> -4592.400000000001
$ws.Cells.Item(136, 8).Value = 4337.3335  # H136: 5750 -> 4337.3335
$ws.Cells.Item(136, 9).Value = 1506  # I136: 1500 -> 1506
$ws.Cells.Item(136, 11).Value = 4518  # K136: 4500 -> 4518
$ws.Cells.Item(136, 13).Value = -1968  # M136: -1950 -> -1968
$ws.Cells.Item(141, 8).Value = 545000  # H141: 393113.66 -> 545000
$ws.Cells.Item(141, 10).Value = 1000000  # J141: 544670.5 -> 1000000
$ws.Cells.Item(141, 12).Value = 1000000  # L141: 544670.5 -> 1000000
$ws.Cells.Item(141, 14).Value = -1010360  # N141: -555030.5 -> -1010360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 86753.414  # H7: 86756 -> 86753.414
$ws.Cells.Item(7, 9).Value = 142893.14  # I7: 142897.58 -> 142893.14
$ws.Cells.Item(7, 11).Value = 428679.42  # K7: 428692.74 -> 428679.42
$ws.Cells.Item(7, 13).Value = -428567.42  # M7: -428580.74 -> -428567.42
$ws.Cells.Item(11, 8).Value = 123.25  # H11: 234 -> 123.25
$ws.Cells.Item(11, 9).Value = 196.5  # I11: 234 -> 196.5
$ws.Cells.Item(11, 10).Value = 50  # J11: 0 -> 50
$ws.Cells.Item(11, 11).Value = 589.5  # K11: 702 -> 589.5
$ws.Cells.Item(11, 12).Value = 150  # L11: 0 -> 150
$ws.Cells.Item(11, 13).Value = -449.5  # M11: -562 -> -449.5
$ws.Cells.Item(11, 14).Value = -430  # N11: None -> -430
$ws.Cells.Item(92, 8).Value = 360  # H92: 402.0909 -> 360
$ws.Cells.Item(92, 9).Value = 381.7  # I92: 403.77777 -> 381.7
$ws.Cells.Item(92, 10).Value = 305.75  # J92: 394.5 -> 305.75
$ws.Cells.Item(92, 11).Value = 1145.1  # K92: 1211.33331 -> 1145.1
$ws.Cells.Item(92, 12).Value = 917.25  # L92: 1183.5 -> 917.25
$ws.Cells.Item(92, 13).Value = 102.9000000000001  # M92: 36.66669000000002 -> 102.9000000000001
$ws.Cells.Item(92, 14).Value = -3413.25  # N92: -3679.5 -> -3413.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 1181.75  # H9: 742.3333 -> 1181.75
$ws.Cells.Item(9, 10).Value = 2500  # J9: 0 -> 2500
$ws.Cells.Item(9, 12).Value = 2500  # L9: 0 -> 2500
$ws.Cells.Item(9, 14).Value = -2840  # N9: None -> -2840
$ws.Cells.Item(43, 8).Value = 19566.467  # H43: 15092.182 -> 19566.467
$ws.Cells.Item(43, 9).Value = 916.6667  # I43: 788.1429000000001 -> 916.6667
$ws.Cells.Item(43, 10).Value = 31999.666  # J43: 40124.25 -> 31999.666
$ws.Cells.Item(43, 11).Value = 916.6667  # K43: 788.1429000000001 -> 916.6667
$ws.Cells.Item(43, 12).Value = 31999.666  # L43: 40124.25 -> 31999.666
$ws.Cells.Item(43, 13).Value = -765.6667  # M43: -637.1429000000001 -> -765.6667
$ws.Cells.Item(43, 14).Value = -32301.666  # N43: -40426.25 -> -32301.666
$ws.Cells.Item(55, 8).Value = 0  # H55: 9000 -> 0
$ws.Cells.Item(55, 9).Value = 0  # I55: 9000 -> 0
$ws.Cells.Item(55, 11).Value = 0  # K55: 9000 -> 0
$ws.Cells.Item(55, 13).Value = $null  # M55: -8673 -> None
$ws.Cells.Item(132, 8).Value = 3031.875  # H132: 3159.0952 -> 3031.875
$ws.Cells.Item(132, 9).Value = 2903.4  # I132: 3037.8823 -> 2903.4
$ws.Cells.Item(132, 11).Value = 8710.200000000001  # K132: 9113.6469 -> 8710.200000000001
$ws.Cells.Item(132, 13).Value = -6180.200000000001  # M132: -6583.6469 -> -6180.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 10187  # H9: 7124.5 -> 10187
$ws.Cells.Item(9, 9).Value = 1000  # I9: 999.5 -> 1000
$ws.Cells.Item(9, 10).Value = 13249.333  # J9: 13249.5 -> 13249.333
$ws.Cells.Item(9, 11).Value = 1000  # K9: 999.5 -> 1000
$ws.Cells.Item(9, 12).Value = 13249.333  # L9: 13249.5 -> 13249.333
$ws.Cells.Item(9, 13).Value = -776  # M9: -775.5 -> -776
$ws.Cells.Item(9, 14).Value = -13697.333  # N9: -13697.5 -> -13697.333
$ws.Cells.Item(100, 8).Value = 1971.625  # H100: 1997.8 -> 1971.625
$ws.Cells.Item(100, 9).Value = 1979  # I100: 1997.5 -> 1979
$ws.Cells.Item(100, 10).Value = 1949.5  # J100: 1999 -> 1949.5
$ws.Cells.Item(100, 11).Value = 1979  # K100: 1997.5 -> 1979
$ws.Cells.Item(100, 12).Value = 1949.5  # L100: 1999 -> 1949.5
$ws.Cells.Item(100, 13).Value = -1438  # M100: -1456.5 -> -1438
$ws.Cells.Item(100, 14).Value = -3031.5  # N100: -3081 -> -3031.5
$ws.Cells.Item(122, 8).Value = 3629.8333  # H122: 3700 -> 3629.8333
$ws.Cells.Item(122, 9).Value = 3629.8333  # I122: 3700 -> 3629.8333
$ws.Cells.Item(122, 11).Value = 10889.4999  # K122: 11100 -> 10889.4999
$ws.Cells.Item(122, 13).Value = -8439.499899999999  # M122: -8650 -> -8439.499899999999
$ws.Cells.Item(132, 8).Value = 3034.6365  # H132: 3284.75 -> 3034.6365
$ws.Cells.Item(132, 9).Value = 2788.1  # I132: 2968.2856 -> 2788.1
$ws.Cells.Item(132, 11).Value = 8364.299999999999  # K132: 8904.856800000001 -> 8364.299999999999
$ws.Cells.Item(132, 13).Value = -5834.299999999999  # M132: -6374.856800000001 -> -5834.299999999999
$ws.Cells.Item(139, 8).Value = 0  # H139: 74999 -> 0
$ws.Cells.Item(139, 9).Value = 0  # I139: 74999 -> 0
$ws.Cells.Item(139, 11).Value = 0  # K139: 74999 -> 0
$ws.Cells.Item(139, 13).Value = $null  # M139: -69859 -> None

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(43, 8).Value = 0  # H43: 20000 -> 0
$ws.Cells.Item(43, 10).Value = 0  # J43: 20000 -> 0
$ws.Cells.Item(43, 12).Value = $null  # L43: 20000 -> None
$ws.Cells.Item(43, 14).Value = 0  # N43: -20298 -> 0
$ws.Cells.Item(54, 8).Value = 19873.25  # H54: 21352.092 -> 19873.25
$ws.Cells.Item(54, 9).Value = 10000  # I54: 9500 -> 10000
$ws.Cells.Item(54, 10).Value = 29746.5  # J54: 28124.715 -> 29746.5
$ws.Cells.Item(54, 11).Value = 10000  # K54: 9500 -> 10000
$ws.Cells.Item(54, 12).Value = 29746.5  # L54: 28124.715 -> 29746.5
$ws.Cells.Item(54, 13).Value = -9480  # M54: -8980 -> -9480
$ws.Cells.Item(54, 14).Value = -30786.5  # N54: -29164.715 -> -30786.5
$ws.Cells.Item(74, 8).Value = 19498.166  # H74: 19997.8 -> 19498.166
$ws.Cells.Item(74, 9).Value = 20330  # I74: 20329.666 -> 20330
$ws.Cells.Item(74, 10).Value = 18666.334  # J74: 19500 -> 18666.334
$ws.Cells.Item(74, 11).Value = 20330  # K74: 20329.666 -> 20330
$ws.Cells.Item(74, 12).Value = 18666.334  # L74: 19500 -> 18666.334
$ws.Cells.Item(74, 13).Value = -19394  # M74: -19393.666 -> -19394
$ws.Cells.Item(74, 14).Value = -20538.334  # N74: -21372 -> -20538.334
$ws.Cells.Item(77, 8).Value = 19498.166  # H77: 19997.8 -> 19498.166
$ws.Cells.Item(77, 9).Value = 20330  # I77: 20329.666 -> 20330
$ws.Cells.Item(77, 10).Value = 18666.334  # J77: 19500 -> 18666.334
$ws.Cells.Item(77, 11).Value = 60990  # K77: 60988.99800000001 -> 60990
$ws.Cells.Item(77, 12).Value = 55999.00199999999  # L77: 58500 -> 55999.00199999999
$ws.Cells.Item(77, 13).Value = -56310  # M77: -56308.99800000001 -> -56310
$ws.Cells.Item(77, 14).Value = -65359.00199999999  # N77: -67860 -> -65359.00199999999
$ws.Cells.Item(80, 8).Value = 15000  # H80: 14333 -> 15000
$ws.Cells.Item(80, 10).Value = 15000  # J80: 14333 -> 15000
$ws.Cells.Item(80, 12).Value = 15000  # L80: 14333 -> 15000
$ws.Cells.Item(80, 14).Value = -16996  # N80: -16329 -> -16996
$ws.Cells.Item(81, 8).Value = 1988.4615  # H81: 2020.9166 -> 1988.4615
$ws.Cells.Item(81, 9).Value = 1904.1666  # I81: 2020.9166 -> 1904.1666
$ws.Cells.Item(81, 10).Value = 3000  # J81: 0 -> 3000
$ws.Cells.Item(81, 11).Value = 3808.3332  # K81: 4041.8332 -> 3808.3332
$ws.Cells.Item(81, 12).Value = 6000  # L81: 0 -> 6000
$ws.Cells.Item(81, 13).Value = -2747.3332  # M81: -2980.8332 -> -2747.3332
$ws.Cells.Item(81, 14).Value = -8122  # N81: None -> -8122
$ws.Cells.Item(83, 8).Value = 15000  # H83: 14333 -> 15000
$ws.Cells.Item(83, 10).Value = 15000  # J83: 14333 -> 15000
$ws.Cells.Item(83, 12).Value = 45000  # L83: 42999 -> 45000
$ws.Cells.Item(83, 14).Value = -54984  # N83: -52983 -> -54984
$ws.Cells.Item(84, 8).Value = 1988.4615  # H84: 2020.9166 -> 1988.4615
$ws.Cells.Item(84, 9).Value = 1904.1666  # I84: 2020.9166 -> 1904.1666
$ws.Cells.Item(84, 10).Value = 3000  # J84: 0 -> 3000
$ws.Cells.Item(84, 11).Value = 19041.666  # K84: 20209.166 -> 19041.666
$ws.Cells.Item(84, 12).Value = 30000  # L84: 0 -> 30000
$ws.Cells.Item(84, 13).Value = -13737.666  # M84: -14905.166 -> -13737.666
$ws.Cells.Item(84, 14).Value = -40608  # N84: None -> -40608
$ws.Cells.Item(126, 8).Value = 4672.722  # H126: 5589.4614 -> 4672.722
$ws.Cells.Item(126, 9).Value = 4624.9287  # I126: 5922.5557 -> 4624.9287
$ws.Cells.Item(126, 11).Value = 13874.7861  # K126: 17767.6671 -> 13874.7861
$ws.Cells.Item(126, 13).Value = -11404.7861  # M126: -15297.6671 -> -11404.7861
$ws.Cells.Item(132, 8).Value = 4208.75  # H132: 4700.6 -> 4208.75
$ws.Cells.Item(132, 9).Value = 3400  # I132: 3871.5715 -> 3400
$ws.Cells.Item(132, 11).Value = 10200  # K132: 11614.7145 -> 10200
$ws.Cells.Item(132, 13).Value = -7670  # M132: -9084.7145 -> -7670
$ws.Cells.Item(136, 8).Value = 26903.666  # H136: 44298.832 -> 26903.666
$ws.Cells.Item(136, 9).Value = 19190.572  # I136: 39498.5 -> 19190.572
$ws.Cells.Item(136, 11).Value = 57571.716  # K136: 118495.5 -> 57571.716
$ws.Cells.Item(136, 13).Value = -55021.716  # M136: -115945.5 -> -55021.716

